$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column A slightly (16.42578125 -> 15.42578125 stored width units)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666

# Update the simulated/recalculated values in column A
$ws.Range("A1").Value = 0.24739049268230673
$ws.Range("A2").Value = -0.0059999999439845908
$ws.Range("A3").Value = -0.0039999999494435556
$ws.Range("A4").Value = -0.0079999999078417261
$ws.Range("A5").Value = -0.0029999999463408145
$ws.Range("A6").Value = -0.001999999941439512
$ws.Range("A7").Value = -0.0099999998703621529
$ws.Range("A8").Value = -0.0099999998665407652
$ws.Range("A9").Value = -0.0019999999346036468
$ws.Range("A10").Value = -0.001999999931543428
$ws.Range("A11").Value = -0.0029999999223768725
$ws.Range("A12").Value = -0.0034999999168761065
$ws.Range("A13").Value = -0.0034999999140454818
$ws.Range("A14").Value = -0.0079999998736504097
$ws.Range("A15").Value = -0.00099999993600441428
$ws.Range("A16").Value = -0.0019999999267561464
$ws.Range("A17").Value = -0.0019999999267188429
$ws.Range("A18").Value = -0.0039999999086690607
$ws.Range("A19").Value = -0.045771414292053336
$ws.Range("A20").Value = -0.0039999999580384582
$ws.Range("A21").Value = -0.0039999999576574297
$ws.Range("A22").Value = -0.0039999999573279155
$ws.Range("A23").Value = -0.0049999999376044713
$ws.Range("A24").Value = -0.019999999798542056
$ws.Range("A25").Value = -0.019999999795939694
$ws.Range("A26").Value = -0.0024999999355515001
$ws.Range("A27").Value = -0.0024999999333812362
$ws.Range("A28").Value = -0.001999999928367302
$ws.Range("A29").Value = -0.00699999987704647
$ws.Range("A30").Value = -0.059999999400994763
$ws.Range("A31").Value = 0.041472851492432738
$ws.Range("A32").Value = -0.0099999998474178398
$ws.Range("A33").Value = -0.018955321092457567
